$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "max" column (C) entirely - its data is being replaced by a
# numeric prediction value that now lives in column B, and the former D/E
# columns collapse down to C/D.
$ws.Columns.Item(3).Delete()

# Header row
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"

# Data rows: column B gets the new numeric values, column C becomes the
# "f__CAG-313" prediction label, column D keeps "f__CAG-313".
$values = @(
    2673.302126080028,
    2003.706563849601,
    2027.092126810354,
    2147.939955225212,
    1920.15760054469,
    1783.420104854294,
    1978.927328704527
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
    $ws.Cells.Item($row, 3).Value = "f__CAG-313"
    $ws.Cells.Item($row, 4).Value = "f__CAG-313"
}
